$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet right after Sheet1 and rename it
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "Code"
$ws2.Range("B1").Value = "Total_Mention"

# Data rows (Code / Total_Mention) matching the order used on Sheet1
$codes  = @("ACCESS","HISTORY","CAPTIVE","REFER","MODEL","RELATED","HUMAN","COLLAB","CONSER","PERMITS","OPTION","METHODS","MAX")
$values = @(63,58,23,27,27,19,54,27,29,26,10,10,63)

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $codes[$i]
    $ws2.Cells.Item($row, 2).Value = $values[$i]
}

# Size column B to fit the "Total_Mention" header, like the authored workbook.
[void]$ws2.Columns("B:B").AutoFit()

# Selection / active-cell bookkeeping to mirror the authored workbook:
# Sheet1 loses its tabSelected flag and its selection becomes the full A1:A14 column.
[void]$ws1.Range("A1:A14").Select()

# Sheet2 becomes the tab-selected / active sheet, with B15 as the "next" cell.
[void]$ws2.Range("B15").Select()
